# Update the dSF (column F) values for several games to reflect the
# repulled / recalculated data (mean calculation) per commit message.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    8  = 1
    14 = -2
    15 = 0
    20 = -2
    54 = 3
    61 = 2
    63 = 2
    64 = 1
    67 = 3
    71 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
